$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "96.182.56"
Set-TextValue "E2" "  +0.59%  "

# Row 3
Set-TextValue "D3" "3.571.34"
Set-TextValue "E3" "  -0.79%  "

# Row 4
Set-TextValue "E4" "  +0.04%  "

# Row 5
Set-TextValue "D5" "240.62"
Set-TextValue "E5" "  +0.91%  "

# Row 6
Set-TextValue "D6" "653.88"
Set-TextValue "E6" "  -0.14%  "

# Row 7
Set-TextValue "E7" "  +7.17%  "

# Row 8
Set-TextValue "D8" "0.404"
Set-TextValue "E8" "  -0.46%  "

# Row 9
Set-TextValue "E9" "  +0.05%  "

# Row 10
Set-TextValue "D10" "1.04"
Set-TextValue "E10" "  +3.43%  "

# Row 11
Set-TextValue "D11" "3.569.58"
Set-TextValue "E11" "  -0.79%  "

# Row 12
Set-TextValue "D12" "43.02"
Set-TextValue "E12" "  -0.09%  "

# Row 13
Set-TextValue "E13" "  +0.74%  "

# Row 14
Set-TextValue "D14" "6.38"
Set-TextValue "E14" "  +1.01%  "

# Row 15
Set-TextValue "D15" "4.234.33"
Set-TextValue "E15" "  -1.23%  "

# Row 16
Set-TextValue "D16" "96.110.77"
Set-TextValue "E16" "  +0.70%  "

# Row 17
Set-TextValue "E17" "  +1.32%  "

# Row 18
Set-TextValue "D18" "3.558.42"
Set-TextValue "E18" "  -1.08%  "

# Row 19
Set-TextValue "D19" "7.75"
Set-TextValue "E19" "  -2.45%  "

# Row 20
Set-TextValue "D20" "12.56"
Set-TextValue "E20" "  -0.10%  "

# Row 21
Set-TextValue "D21" "17.67"
Set-TextValue "E21" "  -1.96%  "

# Row 22
Set-TextValue "D22" "0.507"
Set-TextValue "E22" "  +3.06%  "

# Row 23
Set-TextValue "B23" "SuiNetwork"
Set-TextValue "C23" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D23" "3.41"
Set-TextValue "E23" "  -5.29%  "

# Row 24
Set-TextValue "B24" "BitcoinCash"
Set-TextValue "C24" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D24" "507.21"
Set-TextValue "E24" "  -0.58%  "

# Row 25
Set-TextValue "D25" "0.0000200"
Set-TextValue "E25" "  +2.14%  "

# Row 26
Set-TextValue "D26" "6.85"
Set-TextValue "E26" "  +3.30%  "

# Row 27
Set-TextValue "D27" "95.97"
Set-TextValue "E27" "  -1.00%  "

# Row 28
Set-TextValue "D28" "12.62"
Set-TextValue "E28" "  -0.72%  "

# Row 29
Set-TextValue "D29" "3.762.74"
Set-TextValue "E29" "  -0.40%  "

# Row 30
Set-TextValue "D30" "0.150"
Set-TextValue "E30" "  +7.42%  "

# Row 31
Set-TextValue "D31" "2.98"
Set-TextValue "E31" "  -6.40%  "

# Row 32
Set-TextValue "D32" "11.38"
Set-TextValue "E32" "  +0.62%  "

# Row 33
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  +0.00%  "

# Row 34
Set-TextValue "E34" "  +3.30%  "

# Row 35
Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  -0.16%  "

# Row 36
Set-TextValue "D36" "31.44"
Set-TextValue "E36" "  -1.42%  "

# Row 37
Set-TextValue "D37" "622.48"
Set-TextValue "E37" "  +8.06%  "

# Row 38
Set-TextValue "D38" "8.74"
Set-TextValue "E38" "  +6.91%  "

# Row 39
Set-TextValue "E39" "  +0.52%  "

# Row 40
Set-TextValue "D40" "1.62"
Set-TextValue "E40" "  +8.34%  "

# Row 42
Set-TextValue "E42" "  -0.07%  "

# Row 43
Set-TextValue "E43" "  -2.30%  "

# Row 44
Set-TextValue "E44" "  +5.64%  "

# Row 45
Set-TextValue "E45" "  +2.10%  "

# Row 46
Set-TextValue "D46" "5.68"
Set-TextValue "E46" "  -1.00%  "

# Row 47
Set-TextValue "D47" "23.49"
Set-TextValue "E47" "  -1.19%  "

# Row 48
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "33.96"
Set-TextValue "E48" "  +0.45%  "

# Row 49
Set-TextValue "B49" "VeChain"
Set-TextValue "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0419"
Set-TextValue "E49" "  +0.60%  "

# Row 50
Set-TextValue "D50" "3.54"
Set-TextValue "E50" "  +2.07%  "

# Row 51
Set-TextValue "D51" "8.15"
Set-TextValue "E51" "  +0.42%  "
